$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Insert two new rows before the old row 10 ("Avoid deforestion..." note) to make
# room for a new highlighted callout explaining the EU EPS currently reuses US data.
$ws.Rows("10:11").Insert()

# New row 10: highlighted (yellow) note spanning columns A:D.
$ws.Range("A10").Value = "The EU EPS currently uses US EPS values."
$ws.Range("A10:D10").Interior.Color = 65535

# New row 11: blank spacer row, bold style (matches the other section-header rows
# such as A9/A25 in this sheet).
$ws.Range("A11").Font.Bold = $true
